# Add exercise 1.12 entry to the hour diary

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 9 data
$ws.Cells.Item(9, 1).Value = 211001
$ws.Cells.Item(9, 2).Value = 30
$ws.Cells.Item(9, 3).Value = 1
$ws.Cells.Item(9, 4).Value = "ex 1.12"

# Update the selection to match the post-edit state (next empty row, column B)
$ws.Range("B10").Select()

# Recalculate so the SUM formula in G2 reflects the new data
$excel.Calculate()
